$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 20:52"

# Re-sort country rows 176-211: several countries swapped positions
# as the case counts were refreshed (alphabetised ties resolved by new counts).
# Only the country-name (column A) cells need to move; columns B-H are
# updated below with the refreshed statistics for their (new) row.
$ws.Range("A177").Value = "Laos"
$ws.Range("A178").Value = "Seychelles"
$ws.Range("A181").Value = "Angola"
$ws.Range("A182").Value = "Liberia"
$ws.Range("A183").Value = "Sudan"
$ws.Range("A184").Value = "Republica del Chad"
$ws.Range("A185").Value = "Suazilandia"
$ws.Range("A186").Value = "San Cristobal y Nieves"
$ws.Range("A187").Value = "Zimbabue"
$ws.Range("A188").Value = "Nepal"
$ws.Range("A189").Value = "Montserrat"
$ws.Range("A190").Value = "Republica de Africa Central"
$ws.Range("A193").Value = "Cabo Verde"
$ws.Range("A194").Value = "Somalia"
$ws.Range("A200").Value = "Belice"
$ws.Range("A201").Value = "Sierra Leona"
$ws.Range("A205").Value = "Anguila"
$ws.Range("A206").Value = "Burundi"
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Papua Nueva Guinea"

# Refresh case statistics for the affected rows
# Row 4
$ws.Range("B4").Value = 301147
$ws.Range("C4").Value = 23986
$ws.Range("E4").Value = 278456
$ws.Range("F4").Value = 7973
$ws.Range("G4").Value = 769
$ws.Range("H4").Value = 8173

# Row 43
$ws.Range("D43").Value = 914
$ws.Range("E43").Value = 771

# Row 66
$ws.Range("E66").Value = 753
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 11

# Row 110
$ws.Range("B110").Value = 162
$ws.Range("C110").Value = 7
$ws.Range("D110").Value = 36
$ws.Range("E110").Value = 125

# Row 118
$ws.Range("B118").Value = 134
$ws.Range("C118").Value = 4
$ws.Range("E118").Value = 103

# Row 181
$ws.Range("C181").Value = 2
$ws.Range("D181").Value = 2
$ws.Range("E181").Value = 6
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 2

# Row 182
$ws.Range("C182").Value = 3
$ws.Range("D182").Value = 3
$ws.Range("G182").Value = 1
$ws.Range("H182").Value = 1

# Row 183
$ws.Range("B183").Value = 10
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 2
$ws.Range("E183").Value = 6
$ws.Range("H183").Value = 2

# Row 184
$ws.Range("C184").Value = 1

# Row 186
$ws.Range("E186").Value = 9
$ws.Range("H186").Value = 0

# Row 187
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("H187").Value = 1

# Row 188
$ws.Range("C188").Value = 3
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 8
$ws.Range("H188").Value = 0

# Row 189
$ws.Range("B189").Value = 9
$ws.Range("E189").Value = 7
$ws.Range("H189").Value = 2

# Row 190
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 8
$ws.Range("H190").Value = 0

# Row 193
$ws.Range("C193").Value = 1
$ws.Range("D193").Value = 0
$ws.Range("H193").Value = 1

# Row 194
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 1
$ws.Range("H194").Value = 0

# Row 200
$ws.Range("C200").Value = 0

# Row 201
$ws.Range("C201").Value = 2

